# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the 227fc64f-... row (row 3) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 05:44:27"
$wsZhCn.Range("H3").Value = "2016-03-18 05:44:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 05:44:29"
$wsDeDe.Range("H3").Value = "2016-03-18 05:44:50"
